$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "65.303.77"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  -2.35%  "

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.337.90"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  -4.15%  "

$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = $origStyle
$ws.Range("E4").Value = "  -0.61%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "184.82"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -6.37%  "

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "534.88"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  -1.72%  "

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.610"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  +1.75%  "

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.333.60"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  -3.89%  "

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  -0.22%  "

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.627"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  -3.21%  "

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "59.89"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  -4.10%  "

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.135"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  -4.22%  "

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0000267"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  +0.50%  "

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.20"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  -4.79%  "

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.854.07"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  -5.30%  "

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.324.63"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  -5.15%  "

$ws.Range("E17").Value = "  -4.65%  "

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "17.90"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  -2.01%  "

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "64.905.21"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  -2.83%  "

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.23"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -3.54%  "

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.970"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  -4.86%  "

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "378.43"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  -2.19%  "

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.86"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  -2.59%  "

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.36"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  -3.37%  "

$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "81.57"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  +0.10%  "

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.90"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  +5.69%  "

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.10"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -0.60%  "

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.72"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  -1.79%  "

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.65"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  -3.19%  "

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.55"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  -1.46%  "

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "29.32"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  -3.91%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "650.57"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  -4.40%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.88"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  -2.62%  "

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.42"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -1.81%  "

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "59.94"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  -4.87%  "

$ws.Range("E36").Value = "  -2.78%  "

$ws.Range("E37").Value = "  +0.00%  "

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.399"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  +0.13%  "

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "36.97"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  -3.54%  "

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0736"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  +9.50%  "

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.995"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  -0.52%  "

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.129"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  -0.86%  "

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.926.40"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  -4.41%  "

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.57"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  +2.91%  "

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.73"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  -8.41%  "

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0407"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  +3.77%  "

$ws.Range("E47").Value = "  +12.61%  "

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.68"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -1.31%  "

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.65"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  -1.88%  "

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.128"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  +1.43%  "

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.99"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  +4.08%  "
